$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 5327.615
$ws.Range("I41").Value = 1315
$ws.Range("J41").Value = 7111
$ws.Range("K41").Value = 1315
$ws.Range("L41").Value = 7111
$ws.Range("M41").Value = -875
$ws.Range("N41").Value = -7991

$ws.Range("H121").Value = 776.53845
$ws.Range("J121").Value = 826.8182
$ws.Range("L121").Value = 2480.4546
$ws.Range("N121").Value = -5974.4546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 26316408
$ws.Range("J2").Value = 917.82355
$ws.Range("L2").Value = 917.82355
$ws.Range("N2").Value = -1143.82355

$ws.Range("H32").Value = 18128.55
$ws.Range("I32").Value = 18814.715
$ws.Range("J32").Value = 15246.667
$ws.Range("K32").Value = 18814.715
$ws.Range("L32").Value = 15246.667
$ws.Range("M32").Value = -18527.715
$ws.Range("N32").Value = -15820.667

$ws.Range("H116").Value = 26316408
$ws.Range("J116").Value = 917.82355
$ws.Range("L116").Value = 917.82355
$ws.Range("N116").Value = -5505.82355

$ws.Range("H132").Value = 3635.5535
$ws.Range("I132").Value = 3769.5715
$ws.Range("J132").Value = 3233.5
$ws.Range("K132").Value = 11308.7145
$ws.Range("L132").Value = 9700.5
$ws.Range("M132").Value = -8778.7145
$ws.Range("N132").Value = -14760.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 26316408
$ws.Range("J3").Value = 917.82355
$ws.Range("L3").Value = 917.82355
$ws.Range("N3").Value = -1145.82355

$ws.Range("H20").Value = 2820.0278
$ws.Range("I20").Value = 3345.0435
$ws.Range("J20").Value = 1891.1538
$ws.Range("K20").Value = 3345.0435
$ws.Range("L20").Value = 1891.1538
$ws.Range("M20").Value = -3098.0435
$ws.Range("N20").Value = -2385.1538

$ws.Range("H86").Value = 1853.5
$ws.Range("I86").Value = 1724.2
$ws.Range("J86").Value = 2069
$ws.Range("K86").Value = 1724.2
$ws.Range("L86").Value = 2069
$ws.Range("M86").Value = -601.2
$ws.Range("N86").Value = -4315

$ws.Range("H89").Value = 1853.5
$ws.Range("I89").Value = 1724.2
$ws.Range("J89").Value = 2069
$ws.Range("K89").Value = 8621
$ws.Range("L89").Value = 10345
$ws.Range("M89").Value = -3005
$ws.Range("N89").Value = -21577

$ws.Range("H94").Value = 758.3333
$ws.Range("I94").Value = 1350
$ws.Range("J94").Value = 462.5
$ws.Range("K94").Value = 1350
$ws.Range("L94").Value = 462.5
$ws.Range("M94").Value = -899
$ws.Range("N94").Value = -1364.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 922
$ws.Range("I16").Value = 871.8182
$ws.Range("J16").Value = 983.3333
$ws.Range("K16").Value = 871.8182
$ws.Range("L16").Value = 983.3333
$ws.Range("M16").Value = -584.8182
$ws.Range("N16").Value = -1557.3333

$ws.Range("H31").Value = 2780635
$ws.Range("I31").Value = 2106.6487
$ws.Range("J31").Value = 7250441.5
$ws.Range("K31").Value = 2106.6487
$ws.Range("L31").Value = 7250441.5
$ws.Range("M31").Value = -1811.6487
$ws.Range("N31").Value = -7251031.5

$ws.Range("H34").Value = 2780635
$ws.Range("I34").Value = 2106.6487
$ws.Range("J34").Value = 7250441.5
$ws.Range("K34").Value = 2106.6487
$ws.Range("L34").Value = 7250441.5
$ws.Range("M34").Value = -1904.6487
$ws.Range("N34").Value = -7250845.5

$ws.Range("H113").Value = 922
$ws.Range("I113").Value = 871.8182
$ws.Range("J113").Value = 983.3333
$ws.Range("K113").Value = 871.8182
$ws.Range("L113").Value = 983.3333
$ws.Range("M113").Value = 1298.1818
$ws.Range("N113").Value = -5323.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1188
$ws.Range("I5").Value = 312.33334
$ws.Range("J5").Value = 1592.1538
$ws.Range("K5").Value = 937.0000200000001
$ws.Range("L5").Value = 4776.4614
$ws.Range("M5").Value = -825.0000200000001
$ws.Range("N5").Value = -5000.4614

$ws.Range("H121").Value = 20871.51
$ws.Range("I121").Value = 8601.416999999999
$ws.Range("J121").Value = 24295.72
$ws.Range("K121").Value = 25804.251
$ws.Range("L121").Value = 72887.16
$ws.Range("M121").Value = -24494.251
$ws.Range("N121").Value = -75507.16

$ws.Range("H131").Value = 789.82294
$ws.Range("I131").Value = 573.3333
$ws.Range("J131").Value = 796.80646
$ws.Range("K131").Value = 1719.9999
$ws.Range("L131").Value = 2390.41938
$ws.Range("M131").Value = 3320.0001
$ws.Range("N131").Value = -12470.41938

$ws.Range("H135").Value = 1188
$ws.Range("I135").Value = 312.33334
$ws.Range("J135").Value = 1592.1538
$ws.Range("K135").Value = 2811.00006
$ws.Range("L135").Value = 14329.3842
$ws.Range("M135").Value = -276.0000600000003
$ws.Range("N135").Value = -19399.3842

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 27693.975
$ws.Range("I132").Value = 41683.16
$ws.Range("J132").Value = 2713.2856
$ws.Range("K132").Value = 125049.48
$ws.Range("L132").Value = 8139.8568
$ws.Range("M132").Value = -122519.48
$ws.Range("N132").Value = -13199.8568

$ws.Range("H133").Value = 32000
$ws.Range("J133").Value = 32000
$ws.Range("L133").Value = 32000
$ws.Range("N133").Value = -42120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 523.45
$ws.Range("I22").Value = 575.8182
$ws.Range("J22").Value = 459.44446
$ws.Range("K22").Value = 575.8182
$ws.Range("L22").Value = 459.44446
$ws.Range("M22").Value = -280.8182
$ws.Range("N22").Value = -1049.44446

$ws.Range("H27").Value = 523.45
$ws.Range("I27").Value = 575.8182
$ws.Range("J27").Value = 459.44446
$ws.Range("K27").Value = 575.8182
$ws.Range("L27").Value = 459.44446
$ws.Range("M27").Value = -468.8182
$ws.Range("N27").Value = -673.4444599999999

$ws.Range("H136").Value = 7734.091
$ws.Range("I136").Value = 20751.334
$ws.Range("J136").Value = 2852.625
$ws.Range("K136").Value = 62254.00199999999
$ws.Range("L136").Value = 8557.875
$ws.Range("M136").Value = -59704.00199999999
$ws.Range("N136").Value = -13657.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 15000
$ws.Range("J109").Value = 15000
$ws.Range("L109").Value = 15000
$ws.Range("N109").Value = -17774

$ws.Range("H136").Value = 4032.0264
$ws.Range("I136").Value = 4046.4055
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 12139.2165
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -9589.216499999999
$ws.Range("N136").Value = -15600
